$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells: Wins, Losses, Ties (columns AD, AE, AF) matching the
# existing header formatting (bold, centered, thin border, top aligned)
$headerRange = $ws.Range("AD1:AF1")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows 2-47: Wins=76, Losses=86, Ties=0 for every team row
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
